$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 15.33611787836899
$ws.Range("C2").Value = 9.62931669507082
$ws.Range("D2").Value = 4.597420465463608
$ws.Range("E2").Value = 11.29807201308701
$ws.Range("F2").Value = 23.94278361640686
$ws.Range("I2").Value = 21.21371602361772
$ws.Range("L2").Value = 9.793511073609263
$ws.Range("N2").Value = 16.97128367727474
$ws.Range("O2").Value = 21.30764794258969
$ws.Range("B3").Value = 14.79978730627835
$ws.Range("C3").Value = 9.385398829475919
$ws.Range("D3").Value = 4.561524701291841
$ws.Range("E3").Value = 11.33472574819229
$ws.Range("F3").Value = 23.90313876310368
$ws.Range("I3").Value = 21.30434765447411
$ws.Range("L3").Value = 9.766638514954566
$ws.Range("N3").Value = 17.01751924943386
$ws.Range("O3").Value = 21.33763568331365
$ws.Range("B4").Value = 14.46219974906846
$ws.Range("C4").Value = 9.231070339598851
$ws.Range("D4").Value = 4.539054447938458
$ws.Range("E4").Value = 11.35908074061769
$ws.Range("F4").Value = 23.88630220175017
$ws.Range("I4").Value = 21.36494592026638
$ws.Range("L4").Value = 9.75186690558504
$ws.Range("N4").Value = 17.04769044553765
$ws.Range("O4").Value = 21.36183839125057
$ws.Range("B5").Value = 14.32275817636329
$ws.Range("C5").Value = 9.167091483399329
$ws.Range("D5").Value = 4.529793137811152
$ws.Range("E5").Value = 11.36947058010747
$ws.Range("F5").Value = 23.88133300709764
$ws.Range("I5").Value = 21.39088183560243
$ws.Range("L5").Value = 9.746286305178936
$ws.Range("N5").Value = 17.06043444482736
$ws.Range("O5").Value = 21.37315355256291
$ws.Range("B6").Value = 14.29949746584862
$ws.Range("C6").Value = 9.156403872213788
$ws.Range("D6").Value = 4.528249100109575
$ws.Range("E6").Value = 11.37122388656356
$ws.Range("F6").Value = 23.88062223020153
$ws.Range("I6").Value = 21.39526337482694
$ws.Range("L6").Value = 9.745386281979393
$ws.Range("N6").Value = 17.06257771805997
$ws.Range("O6").Value = 21.37512002751485
$ws.Range("B7").Value = 14.4603264789444
$ws.Range("C7").Value = 9.230211826494209
$ws.Range("D7").Value = 4.538929964947791
$ws.Range("E7").Value = 11.35921897902846
$ws.Range("F7").Value = 23.88622752105204
$ws.Range("I7").Value = 21.36529067791629
$ws.Range("L7").Value = 9.751789860890664
$ws.Range("N7").Value = 17.04786049654992
$ws.Range("O7").Value = 21.36198511618566
$ws.Range("B8").Value = 15.15302908782256
$ws.Range("C8").Value = 9.546194401480383
$ws.Range("D8").Value = 4.585135055784989
$ws.Range("E8").Value = 11.31032632289195
$ws.Range("F8").Value = 23.92755961959307
$ws.Range("I8").Value = 21.24393656076753
$ws.Range("L8").Value = 9.783889455423596
$ws.Range("N8").Value = 16.98685625101474
$ws.Range("O8").Value = 21.31678399483671
$ws.Range("B9").Value = 16.43750218531817
$ws.Range("C9").Value = 10.12728418958737
$ws.Range("D9").Value = 4.672172816317947
$ws.Range("E9").Value = 11.22912585648783
$ws.Range("F9").Value = 24.06787303428634
$ws.Range("I9").Value = 21.04537618889175
$ws.Range("L9").Value = 9.860340714839054
$ws.Range("N9").Value = 16.88133364027148
$ws.Range("O9").Value = 21.27422101782892
$ws.Range("B10").Value = 17.32645871805533
$ws.Range("C10").Value = 10.52785043375107
$ws.Range("D10").Value = 4.73373651895815
$ws.Range("E10").Value = 11.17842023797238
$ws.Range("F10").Value = 24.20656221360428
$ws.Range("I10").Value = 20.92370960583138
$ws.Range("L10").Value = 9.924438123534463
$ws.Range("N10").Value = 16.81235697984877
$ws.Range("O10").Value = 21.27118249664104
$ws.Range("B11").Value = 17.71729242500076
$ws.Range("C11").Value = 10.70383091105485
$ws.Range("D11").Value = 4.761183186996946
$ws.Range("E11").Value = 11.15729778451201
$ws.Range("F11").Value = 24.27723070043848
$ws.Range("I11").Value = 20.87365832169452
$ws.Range("L11").Value = 9.955246902742152
$ws.Range("N11").Value = 16.78282440021101
$ws.Range("O11").Value = 21.27594640141345
$ws.Range("B12").Value = 17.86322277005385
$ws.Range("C12").Value = 10.76953435833005
$ws.Range("D12").Value = 4.771492574080225
$ws.Range("E12").Value = 11.14957887830411
$ws.Range("F12").Value = 24.30506448100369
$ws.Range("I12").Value = 20.85546997235775
$ws.Range("L12").Value = 9.967144042891682
$ws.Range("N12").Value = 16.77190580671756
$ws.Range("O12").Value = 21.27863415369143
$ws.Range("B13").Value = 17.83188787783763
$ws.Range("C13").Value = 10.75542618570655
$ws.Range("D13").Value = 4.76927606019276
$ws.Range("E13").Value = 11.15122883972077
$ws.Range("F13").Value = 24.29902253355199
$ws.Range("I13").Value = 20.85935308049415
$ws.Range("L13").Value = 9.964571640348311
$ws.Range("N13").Value = 16.77424555538198
$ws.Range("O13").Value = 21.27801600153713
$ws.Range("B14").Value = 17.72934032291318
$ws.Range("C14").Value = 10.70925535517284
$ws.Range("D14").Value = 4.762033052267792
$ws.Range("E14").Value = 11.15665713858387
$ws.Range("F14").Value = 24.27949919635165
$ws.Range("I14").Value = 20.8721466038634
$ws.Range("L14").Value = 9.956221107807236
$ws.Range("N14").Value = 16.78192081820759
$ws.Range("O14").Value = 21.27614981549453
$ws.Range("B15").Value = 17.6662541041747
$ws.Range("C15").Value = 10.68085131821112
$ws.Range("D15").Value = 4.757585443963308
$ws.Range("E15").Value = 11.16001856035465
$ws.Range("F15").Value = 24.26767982094022
$ws.Range("I15").Value = 20.88008273339046
$ws.Range("L15").Value = 9.951135978514651
$ws.Range("N15").Value = 16.78665659985344
$ws.Range("O15").Value = 21.27512180021426
$ws.Range("B16").Value = 17.30063337584057
$ws.Range("C16").Value = 10.51622069868681
$ws.Range("D16").Value = 4.731931249644134
$ws.Range("E16").Value = 11.17983977067527
$ws.Range("F16").Value = 24.20209498101462
$ws.Range("I16").Value = 20.92708743365626
$ws.Range("L16").Value = 9.922457334884111
$ws.Range("N16").Value = 16.8143240834834
$ws.Range("O16").Value = 21.2709948546619
$ws.Range("B17").Value = 17.07277410579498
$ws.Range("C17").Value = 10.41359796165364
$ws.Range("D17").Value = 4.716047418964769
$ws.Range("E17").Value = 11.19249743155485
$ws.Range("F17").Value = 24.16379057342125
$ws.Range("I17").Value = 20.95728209487194
$ws.Range("L17").Value = 9.905281911648629
$ws.Range("N17").Value = 16.83176937360361
$ws.Range("O17").Value = 21.27003750075432
$ws.Range("B18").Value = 16.94044561705245
$ws.Range("C18").Value = 10.35398693253343
$ws.Range("D18").Value = 4.706859063473753
$ws.Range("E18").Value = 11.19996073042437
$ws.Range("F18").Value = 24.14247329386921
$ws.Range("I18").Value = 20.97514749131805
$ws.Range("L18").Value = 9.895558849623532
$ws.Range("N18").Value = 16.84197716057628
$ws.Range("O18").Value = 21.27006541231152
$ws.Range("B19").Value = 16.89542728139908
$ws.Range("C19").Value = 10.33370442376325
$ws.Range("D19").Value = 4.703739158687688
$ws.Range("E19").Value = 11.20251909051802
$ws.Range("F19").Value = 24.13537880126172
$ws.Range("I19").Value = 20.98128187091968
$ws.Range("L19").Value = 9.892293746940371
$ws.Range("N19").Value = 16.8454631963813
$ws.Range("O19").Value = 21.270174218738
$ws.Range("B20").Value = 17.09716246000513
$ws.Range("C20").Value = 10.42458318446018
$ws.Range("D20").Value = 4.717743729618228
$ws.Range("E20").Value = 11.19113106618417
$ws.Range("F20").Value = 24.16779431878116
$ws.Range("I20").Value = 20.95401623158053
$ws.Range("L20").Value = 9.907094189768468
$ws.Range("N20").Value = 16.82989431792153
$ws.Range("O20").Value = 21.27007953210808
$ws.Range("B21").Value = 17.75951808778803
$ws.Range("C21").Value = 10.72284256007626
$ws.Range("D21").Value = 4.76416280998507
$ws.Range("E21").Value = 11.15505512289455
$ws.Range("F21").Value = 24.28520469008918
$ws.Range("I21").Value = 20.86836804670433
$ws.Range("L21").Value = 9.958667661257149
$ws.Range("N21").Value = 16.77965922673338
$ws.Range("O21").Value = 21.27667397869557
$ws.Range("B22").Value = 18.18029522282185
$ws.Range("C22").Value = 10.9122977456256
$ws.Range("D22").Value = 4.794008674159678
$ws.Range("E22").Value = 11.13310784980681
$ws.Range("F22").Value = 24.36818531597045
$ws.Range("I22").Value = 20.81685252482834
$ws.Range("L22").Value = 9.993714188413758
$ws.Range("N22").Value = 16.74837061769356
$ws.Range("O22").Value = 21.28613451076076
$ws.Range("B23").Value = 17.95686265291783
$ws.Range("C23").Value = 10.81169492602865
$ws.Range("D23").Value = 4.778125565317975
$ws.Range("E23").Value = 11.1446722742522
$ws.Range("F23").Value = 24.32333143619523
$ws.Range("I23").Value = 20.84393804348627
$ws.Range("L23").Value = 9.974888897993219
$ws.Range("N23").Value = 16.76492894891686
$ws.Range("O23").Value = 21.28061417753123
$ws.Range("B24").Value = 17.08614061578533
$ws.Range("C24").Value = 10.41961866802458
$ws.Range("D24").Value = 4.716977003362753
$ws.Range("E24").Value = 11.19174821986418
$ws.Range("F24").Value = 24.16598202998836
$ws.Range("I24").Value = 20.95549115219566
$ws.Range("L24").Value = 9.906274386982664
$ws.Range("N24").Value = 16.83074147527696
$ws.Range("O24").Value = 21.27005872832712
$ws.Range("B25").Value = 16.09900253995024
$ws.Range("C25").Value = 9.974522196392904
$ws.Range("D25").Value = 4.649030497439746
$ws.Range("E25").Value = 11.24952090739074
$ws.Range("F25").Value = 24.02361631517392
$ws.Range("I25").Value = 21.09485296442495
$ws.Range("L25").Value = 9.838242758349752
$ws.Range("N25").Value = 16.90837524040976
$ws.Range("O25").Value = 21.28078502532818
